# ------------------------------------------------------------------
# counter_ver.1.xlsx edit: add a "print" sheet, switch the date column
# (B) from literal/custom-formatted text to real dates, append new
# rows 4-12 of data, and move the trailing "**" marker down to row 13.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("2017")

# --- 1. Re-key column B on row 1: a genuine Excel date -------------
# serial 43067 == 2017-11-28
$ws.Cells.Item(1, 2).Value = 43067
$ws.Cells.Item(1, 2).NumberFormat = "mm-dd-yy"

# Copy that date formatting onto B2:B10 so every date cell in the
# column shares one style (mirrors the column's own style=1 default).
$ws.Cells.Item(1, 2).Copy()
$ws.Range("B2:B10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 2 and 3 keep their literal text "2017/11/28 " (trailing space)
$ws.Cells.Item(2, 2).Value = "2017/11/28 "
$ws.Cells.Item(3, 2).Value = "2017/11/28 "

# --- 2. Move the old trailing "**" cell out of row 4 ---------------
$ws.Cells.Item(4, 1).ClearContents()

# --- 3. New data rows 4-10 ------------------------------------------
$ws.Cells.Item(4, 1).Value  = 80001234
$ws.Cells.Item(4, 2).Value  = "2017/11/28 "
$ws.Cells.Item(4, 3).Value  = 2
$ws.Cells.Item(4, 4).Value  = 56
$ws.Cells.Item(4, 5).Value  = 200

$ws.Cells.Item(5, 1).Value  = 80001258
$ws.Cells.Item(5, 2).Value  = "2017/11/28 "
$ws.Cells.Item(5, 3).Value  = 1
$ws.Cells.Item(5, 4).Value  = 23
$ws.Cells.Item(5, 5).Value  = 56

$ws.Cells.Item(6, 1).Value  = 80009999
$ws.Cells.Item(6, 2).Value  = 43067
$ws.Cells.Item(6, 3).Value  = 1
$ws.Cells.Item(6, 4).Value  = 56
$ws.Cells.Item(6, 5).Value  = 60

$ws.Cells.Item(7, 1).Value  = 80005555
$ws.Cells.Item(7, 2).Value  = 43067
$ws.Cells.Item(7, 3).Value  = 6
$ws.Cells.Item(7, 4).Value  = 6
$ws.Cells.Item(7, 5).Value  = 12

$ws.Cells.Item(8, 1).Value  = 89
$ws.Cells.Item(8, 2).Value  = 43067
$ws.Cells.Item(8, 3).Value  = 89
$ws.Cells.Item(8, 4).Value  = 89
$ws.Cells.Item(8, 5).Value  = 178

$ws.Cells.Item(9, 1).Value  = 1
$ws.Cells.Item(9, 2).Value  = 43067
$ws.Cells.Item(9, 3).Value  = 1
$ws.Cells.Item(9, 4).Value  = 1
$ws.Cells.Item(9, 5).Value  = 2

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 43067
$ws.Cells.Item(10, 3).Value = 8
$ws.Cells.Item(10, 4).Value = 8
$ws.Cells.Item(10, 5).Value = 16

# Rows 11-12: plain text "2017/11/28" (no trailing space). Enter it as
# a formula returning the literal string, then freeze the formula down
# to a static value via copy/paste-values -- this keeps Excel's
# autodetect from turning the text into a serial date, without
# stamping a brand-new number style.
$ws.Cells.Item(11, 2).Formula = "=""2017/11/28"""
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 2

$ws.Cells.Item(12, 2).Formula = "=""2017/11/28"""
$ws.Cells.Item(12, 2).Copy()
$ws.Cells.Item(12, 2).PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 3).Value = 8
$ws.Cells.Item(12, 4).Value = 8
$ws.Cells.Item(12, 5).Value = 16

# Row 13: the "**" marker that used to sit in row 4.
$ws.Cells.Item(13, 1).Value = "**"

# --- 4. Selection on the "2017" sheet lands on B10 ------------------
$ws.Range("B10").Select()

# --- 5. Add the "print" sheet after "2017" and make it active -------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "print"

$ws2017 = $wb.Worksheets.Item("2017")
$newSheet.Move($null, $ws2017)

$wsPrint = $wb.Worksheets.Item("print")
$wsPrint.Range("J20").Select()
$wsPrint.Activate()
